# Adding GESS to the forecast portfolio
# Shift the timestamp column (A) forward by 9 days and zero out the
# notified production values (B) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value2 = $cellA.Value2 + 9

    $ws.Cells.Item($r, 2).Value2 = 0
}
